$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    Write-Host "Replace '$find' -> '$replace': $result"
}

Replace-Text "08/2017 – 06/2021" "08/2020 – 05/2024"
Replace-Text "Clark University ’21, Worcester, MA, USA" "Clark University ’24, Worcester, MA, USA"
Replace-Text "Dean’s List: Fall 2017 – Spring 2020" "Dean’s List: Fall 2021 – Spring 2023"
Replace-Text "05/2019 " "05/2022 "
Replace-Text "06/2020 – 08/2020" "06/2023 – 08/2023"
Replace-Text "06/2019 – 08/2019" "06/2022 – 08/2022"
Replace-Text "10/2018 – 05/2020" "10/2021 – 05/2022"
Replace-Text "08/2018 – 05/2019" "08/2021 – 05/2022"
Replace-Text "01/2018 – present" "01/2021 – present"
Replace-Text "10/2017 – 10/2019" "10/2020 – 10/2022"
